# Rename/refresh the HNF status rows: turn the raw date-serial values in
# column E into plain text timestamps, reset D4's flag back to 0, and drop
# the now-unused custom date style from E2/E4 so they fall back to the
# default (unstyled) cell format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2: was a styled numeric date serial -> becomes a plain text timestamp.
$ws.Range("E2").Style = "Normal"
$ws.Range("E2").Value = "2022-06-21 09:30:54"

# E3: was an unstyled numeric value -> becomes a plain text timestamp.
$ws.Range("E3").Value = "2022-06-17 11:00:06"

# D4: flag flips from 1 back to 0.
$ws.Range("D4").Value = 0

# E4: was a styled numeric date serial -> becomes a plain text timestamp.
$ws.Range("E4").Style = "Normal"
$ws.Range("E4").Value = "2022-06-19 12:00:08"
